# Progress Diary - fill in Lisa Bott's column (column C) for the "Week 2" block.
# Before this edit, column C (rows 5-9) still held placeholder / unrelated filler
# text ("Name", "Research macht Spass", a phone-call note, "30", and an impediment
# about phone calls). This change replaces that placeholder content with Lisa
# Bott's actual progress-diary entries, matching the layout already used for the
# other two team members in columns B (Gloria Bichler) and D (Marcia Perez).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name of the third team member
$ws.Range("C5").Value = "Lisa Bott"

# Was ich gelernt habe
$ws.Range("C6").Value = "Wie man die Persönlichkeit des Chatbots und/oder Nutzer an das Branding anpassen kann und Designentschei-dungen abwägen muss"

# Mein persönlicher Beitrag
$ws.Range("C7").Value = "Angeschaut woraus Branding aufgebaut ist, Keywords für Brand-Personality aufgeschrieben & Farb-schema, Logo, Icon, Hintergrund/Muster, Schriftarten herausgesucht/ erstellt und ein Mockup für die Implementierung der Design-Vorlage gemacht. Story für zwei Storyboards mit Personas ausgedacht und erstellt."

# Mein zeitlicher Aufwand (in Minuten)
$ws.Range("C8").Value = 450

# Meine Impediments
$ws.Range("C9").Value = "Eine gute Zeiteinteilung kann sehr viel Stress reduzieren."

# Keep column C's formatting consistent with the rest of the "Week 2" block
# (matches the wrap/valign styling already used by columns B and D in this row
# band) so the new text wraps and top-aligns the same way.
$ws.Range("C5:C9").WrapText = $true
$ws.Range("C5:C9").VerticalAlignment = -4160  # xlTop
$ws.Range("C5").Font.Size = 10
$ws.Range("C5").Font.Name = "Arial"
$ws.Range("C6").Font.Size = 10
$ws.Range("C6").Font.Name = "Arial"
$ws.Range("C7").Font.Size = 10
$ws.Range("C7").Font.Name = "Arial"
$ws.Range("C9").Font.Size = 10
$ws.Range("C9").Font.Name = "Arial"

# The longer diary text now in rows 7 and 9 needs a bit more/less vertical
# room than the placeholder text it replaced, so adjust those row heights.
$ws.Rows.Item(7).RowHeight = 182
$ws.Rows.Item(9).RowHeight = 42

# Update the worksheet selection / zoom to reflect the state captured after the edit
$ws.Range("E11").Select()
$ws.Application.ActiveWindow.Zoom = 150
